$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1 - copy format from neighboring header cell G1, then set text
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# New data column H2:H5 ("Save" indicator values)
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 1
